$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / link / percentage updates (safe as literal text)
$ws.Range("D2").Value = '61.219.97'
$ws.Range("E2").Value = '  +0.99%  '
$ws.Range("D3").Value = '2.670.63'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E5").Value = '  +4.01%  '
$ws.Range("E6").Value = '  +1.76%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E9").Value = '  -2.61%  '
$ws.Range("E10").Value = '  +5.11%  '
$ws.Range("E11").Value = '  +2.34%  '
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").Value = '3.128.33'
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").Value = '61.239.03'
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("E15").Value = '  +2.57%  '
$ws.Range("E16").Value = '  +2.31%  '
$ws.Range("D17").Value = '2.676.26'
$ws.Range("E17").Value = '  +2.16%  '
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("E20").Value = '  +1.63%  '
$ws.Range("E21").Value = '  +2.57%  '
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("E24").Value = '  +2.68%  '
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").Value = '0.0₃0869'
$ws.Range("E27").Value = '  +3.07%  '
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("E30").Value = '  +6.47%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E32").Value = '  +4.09%  '
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("E34").Value = '  +4.65%  '
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("E36").Value = '  +8.91%  '
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("E38").Value = '  +6.28%  '
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("E41").Value = '  +4.23%  '
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("E44").Value = '  +2.51%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E45").Value = '  +3.15%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("E47").Value = '  +2.82%  '
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("E49").Value = '  +8.75%  '
$ws.Range("D50").Value = '2.001.25'
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("E51").Value = '  +3.16%  '

# Numeric-looking price text must be forced to remain text,
# otherwise Excel auto-converts it to a floating point number
# and loses formatting (e.g. trailing zeros, multi-dot grouping).
# NumberFormat/ClearFormats are applied per-cell (union Range() args
# only affect the first area in this automation layer).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D5").Value = '533.04'
$ws.Range("D6").Value = '157.09'
$ws.Range("D7").Value = '0.996'
$ws.Range("D8").Value = '0.592'
$ws.Range("D9").Value = '6.62'
$ws.Range("D10").Value = '0.110'
$ws.Range("D15").Value = '22.19'
$ws.Range("D18").Value = '4.80'
$ws.Range("D19").Value = '357.50'
$ws.Range("D20").Value = '10.76'
$ws.Range("D21").Value = '6.34'
$ws.Range("D22").Value = '1.00'
$ws.Range("D23").Value = '61.68'
$ws.Range("D24").Value = '0.435'
$ws.Range("D25").Value = '0.170'
$ws.Range("D28").Value = '7.45'
$ws.Range("D30").Value = '6.20'
$ws.Range("D31").Value = '19.64'
$ws.Range("D32").Value = '1.64'
$ws.Range("D33").Value = '150.37'
$ws.Range("D34").Value = '4.17'
$ws.Range("D35").Value = '1.21'
$ws.Range("D36").Value = '0.919'
$ws.Range("D37").Value = '0.885'
$ws.Range("D38").Value = '309.83'
$ws.Range("D40").Value = '3.83'
$ws.Range("D41").Value = '0.652'
$ws.Range("D43").Value = '20.64'
$ws.Range("D44").Value = '0.0568'
$ws.Range("D45").Value = '5.06'
$ws.Range("D46").Value = '0.997'
$ws.Range("D48").Value = '10.37'
$ws.Range("D49").Value = '19.15'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()

